# docs: update calculator_documentation.docx for feat: add factorial operation
$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Append a bullet line about the new Factorial operation to the
#    "4. Feature Specifications" table-of-contents entry (List Number style).
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $txt = $p.Range.Text
    if ($txt -like "*4. Feature Specifications*") {
        if ($p.Style.NameLocal -eq "List Number") {
            $r = $p.Range
            [void]$r.MoveEnd(1, -1)    # exclude the paragraph mark
            $r.Collapse(0)             # collapse to the end of the text
            $r.InsertAfter([char]11)   # manual line break (w:br)
            $r.Collapse(0)
            $r.InsertAfter([char]0x2022 + " Factorial Operation: Calculates the factorial of a number")
        }
    }
}

# ---------------------------------------------------------------------------
# 2. Add a new "Factorial" / "factorial" row to the Menu Structure table
#    (the table whose last row reads "8" | "Exit").
# ---------------------------------------------------------------------------
foreach ($t in $d.Tables) {
    $lastRow = $t.Rows.Item($t.Rows.Count)
    $lastCellText = $lastRow.Cells.Item($lastRow.Cells.Count).Range.Text
    if ($lastCellText -like "Exit*") {
        $newRow = $t.Rows.Add()
        $newRow.Cells.Item(1).Range.Text = "Factorial"
        $newRow.Cells.Item(2).Range.Text = "factorial"
    }
}
